# Auto update stock data: bump the report date from 2025/12/26 to 2025/12/27
# and refresh the EBITDA readings for the rows whose values moved.
#
# The Date_1 and EBITDA columns are stored as text (quote-prefixed) in the
# workbook, so values are written with a leading apostrophe to keep them as
# text instead of letting Excel auto-convert date-looking / numeric-looking
# strings into a date serial or a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = "2025/12/27"

# Rows whose EBITDA (column B) value also changed.
$ebitdaUpdates = @{
    2  = "6.75"
    20 = "13.00"
    26 = "11.30"
    32 = "27.86"
    44 = "11.20"
    50 = "11.57"
    62 = "11.76"
    68 = "13.17"
    74 = "16.76"
}

# Every data row touched by the refresh (date always updates).
$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($row in $rows) {
    $ws.Range("A$row").Value = "'" + $newDate

    if ($ebitdaUpdates.ContainsKey($row)) {
        $ws.Range("B$row").Value = "'" + $ebitdaUpdates[$row]
    }
}
